$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Date
$ws.Range("B3").Value = Get-Date -Year 2020 -Month 10 -Day 22 -Hour 0 -Minute 0 -Second 0

# Tasks to complete next week (column B, rows 19-20)
$ws.Range("B19").Value = "Create study protocol and questionnaires"
$ws.Range("B20").Value = "Conduct user study with at least 8 users"

# Tasks completed this week (column A, row 19)
$ws.Range("A19").Value = "Picked ""perfect"" idea and prototyped it"

# Team name and number of members
$ws.Range("B4").Value = "Limette"
$ws.Range("B5").Value = 4

# Team member names
$ws.Range("A8").Value = "Lukas Hasler"
$ws.Range("A9").Value = "Pascal Strebel"
$ws.Range("A10").Value = "Cedric Weibel"
$ws.Range("A11").Value = "Robin Schmidiger"

# Team member salaries
$ws.Range("B8").Value = 100
$ws.Range("B9").Value = 100
$ws.Range("B10").Value = 100
$ws.Range("B11").Value = 100

# Clear the old "Member 5" label - only 4 team members now
$ws.Range("A12").Value = ""

# Update selection to match the target workbook
$ws.Range("C13").Select()
